# Updated cryptos list on Sun May 14 15:26:46 UTC 2023 with GitHub Actions
#
# This script refreshes the Price (column D) and Volume(1h) (column E)
# figures for the crypto ranking table, and fixes the relative ordering
# of two pairs of rows (Hedera/VeChain at rows 38-39 and
# NEARProtocol/PaxDollar at rows 48-49) whose rank swapped position.
#
# Column D holds values that look numeric (e.g. "1.035", "27.885.23")
# but must stay plain text exactly as scraped (leading zeros, two-dot
# "thousands" separators, trailing-zero precision, etc. must be kept
# verbatim). Excel's COM layer auto-converts such strings to real
# numbers on assignment, so we temporarily force column D to Text
# number format before writing the values, then restore the original
# "Normal" style/General format afterwards so the cell formatting is
# left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.885.23'
$ws.Range("E2").Value = '  +1.85%  '
$ws.Range("D3").Value = '1.869.06'
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("D4").Value = '1.035'
$ws.Range("E4").Value = '  +1.12%  '
$ws.Range("D5").Value = '324.33'
$ws.Range("E5").Value = '  +1.81%  '
$ws.Range("D6").Value = '1.031'
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("D7").Value = '0.4431'
$ws.Range("E7").Value = '  +1.66%  '
$ws.Range("D8").Value = '0.3823'
$ws.Range("E8").Value = '  +2.71%  '
$ws.Range("D9").Value = '0.07484'
$ws.Range("D10").Value = '0.8913'
$ws.Range("E10").Value = '  +2.30%  '
$ws.Range("D11").Value = '21.80'
$ws.Range("E11").Value = '  +2.37%  '
$ws.Range("D12").Value = '1.888.01'
$ws.Range("E12").Value = '  -2.91%  '
$ws.Range("D13").Value = '5.600'
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("D14").Value = '6.788'
$ws.Range("D15").Value = '0.07216'
$ws.Range("E15").Value = '  +1.24%  '
$ws.Range("D16").Value = '84.89'
$ws.Range("E16").Value = '  +3.39%  '
$ws.Range("D17").Value = '1.035'
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("E18").Value = '  +1.90%  '
$ws.Range("D19").Value = '1.031'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = '15.63'
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("D21").Value = '27.919.08'
$ws.Range("E21").Value = '  +1.86%  '
$ws.Range("D22").Value = '5.340'
$ws.Range("E22").Value = '  +1.91%  '
$ws.Range("D23").Value = '11.36'
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("D24").Value = '2.112.04'
$ws.Range("E24").Value = '  -1.35%  '
$ws.Range("D25").Value = '2.023'
$ws.Range("E25").Value = '  +6.83%  '
$ws.Range("D26").Value = '158.51'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").Value = '18.92'
$ws.Range("E27").Value = '  +2.15%  '
$ws.Range("D28").Value = '5.418'
$ws.Range("E28").Value = '  +3.50%  '
$ws.Range("D29").Value = '1.993'
$ws.Range("E29").Value = '  +4.11%  '
$ws.Range("D30").Value = '118.46'
$ws.Range("E30").Value = '  +2.65%  '
$ws.Range("D31").Value = '0.09074'
$ws.Range("D32").Value = '1.240'
$ws.Range("E32").Value = '  +3.40%  '
$ws.Range("D33").Value = '0.7831'
$ws.Range("E33").Value = '  +3.25%  '
$ws.Range("D34").Value = '4.617'
$ws.Range("E34").Value = '  +3.57%  '
$ws.Range("D35").Value = '3.011'
$ws.Range("E35").Value = '  +5.29%  '
$ws.Range("D36").Value = '1.032'
$ws.Range("D37").Value = '1.150'
$ws.Range("E37").Value = '  +0.08%  '

# Hedera / VeChain swapped ranking position
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01993'
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.05378'
$ws.Range("E39").Value = '  +2.70%  '

$ws.Range("E40").Value = '  +3.84%  '
$ws.Range("D41").Value = '0.5238'
$ws.Range("E41").Value = '  +1.48%  '
$ws.Range("D42").Value = '0.1700'
$ws.Range("E42").Value = '  +2.44%  '
$ws.Range("D43").Value = '6.918'
$ws.Range("E43").Value = '  +6.05%  '
$ws.Range("D44").Value = '8.926'
$ws.Range("E44").Value = '  +5.81%  '
$ws.Range("D45").Value = '112.15'
$ws.Range("E45").Value = '  +3.51%  '
$ws.Range("D46").Value = '10.77'
$ws.Range("E46").Value = '  +2.29%  '
$ws.Range("D47").Value = '0.06624'
$ws.Range("E47").Value = '  +5.31%  '

# NEARProtocol / PaxDollar swapped ranking position
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '1.033'
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.726'
$ws.Range("E49").Value = '  +3.47%  '

$ws.Range("D50").Value = '0.4754'
$ws.Range("E50").Value = '  +2.91%  '
$ws.Range("D51").Value = '1.916'
$ws.Range("E51").Value = '  +2.25%  '

# Restore column D's original "Normal" style / General number format so
# we don't leave a stray text-format style behind on these cells.
$priceRange.Style = "Normal"
